$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2124352331606218
$ws.Range("C2").Value = 0.5336787564766839
$ws.Range("J2").Value = 0.007772020725388601
$ws.Range("P2").Value = 0.1450777202072539
$ws.Range("S2").Value = 0.1010362694300518
$ws.Range("C3").Value = 0.04072398190045249
$ws.Range("J3").Value = 0.04524886877828054
$ws.Range("P3").Value = 0.7285067873303167
$ws.Range("S3").Value = 0.1855203619909502
$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.8108108108108109
$ws.Range("S4").Value = 0.1621621621621622
$ws.Range("B6").Value = 0.06299212598425197
$ws.Range("D6").Value = 0.003937007874015748
$ws.Range("E6").Value = 0.003937007874015748
$ws.Range("F6").Value = 0.05511811023622047
$ws.Range("J6").Value = 0.2283464566929134
$ws.Range("O6").Value = 0.03543307086614173
$ws.Range("Q6").Value = 0.1889763779527559
$ws.Range("R6").Value = 0.05905511811023622
$ws.Range("S6").Value = 0.3622047244094488
$ws.Range("B7").Value = 0.08900523560209424
$ws.Range("D7").Value = 0.01570680628272251
$ws.Range("E7").Value = 0.005235602094240838
$ws.Range("F7").Value = 0.08900523560209424
$ws.Range("J7").Value = 0.1308900523560209
$ws.Range("O7").Value = 0.005235602094240838
$ws.Range("Q7").Value = 0.1518324607329843
$ws.Range("R7").Value = 0.06282722513089005
$ws.Range("S7").Value = 0.450261780104712
$ws.Range("B8").Value = 0.1241050119331742
$ws.Range("D8").Value = 0.01909307875894988
$ws.Range("F8").Value = 0.081145584725537
$ws.Range("J8").Value = 0.08353221957040573
$ws.Range("O8").Value = 0.007159904534606206
$ws.Range("Q8").Value = 0.1909307875894988
$ws.Range("R8").Value = 0.05250596658711217
$ws.Range("S8").Value = 0.441527446300716
$ws.Range("B9").Value = 0.1022222222222222
$ws.Range("D9").Value = 0.02222222222222222
$ws.Range("F9").Value = 0.05333333333333334
$ws.Range("J9").Value = 0.1066666666666667
$ws.Range("O9").Value = 0.01333333333333333
$ws.Range("Q9").Value = 0.1822222222222222
$ws.Range("R9").Value = 0.09333333333333334
$ws.Range("S9").Value = 0.4266666666666667
$ws.Range("B10").Value = 0.1438304314912945
$ws.Range("D10").Value = 0.0174110522331567
$ws.Range("E10").Value = 0.001514004542013626
$ws.Range("F10").Value = 0.07494322482967448
$ws.Range("J10").Value = 0.1097653292959879
$ws.Range("O10").Value = 0.01665404996214989
$ws.Range("Q10").Value = 0.2172596517789553
$ws.Range("R10").Value = 0.07267221801665405
$ws.Range("S10").Value = 0.3459500378501135
$ws.Range("G11").Value = 0.1359773371104816
$ws.Range("J11").Value = 0.1076487252124646
$ws.Range("K11").Value = 0.2096317280453258
$ws.Range("L11").Value = 0.5184135977337111
$ws.Range("S11").Value = 0.028328611898017
$ws.Range("G12").Value = 0.675531914893617
$ws.Range("J12").Value = 0.2287234042553191
$ws.Range("K12").Value = 0.01063829787234043
$ws.Range("L12").Value = 0.02659574468085106
$ws.Range("S12").Value = 0.05851063829787234
$ws.Range("F13").Value = 0.02857142857142857
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.4
$ws.Range("F15").Value = 0.04910714285714286
$ws.Range("H15").Value = 0.1785714285714286
$ws.Range("I15").Value = 0.07142857142857142
$ws.Range("J15").Value = 0.3303571428571428
$ws.Range("K15").Value = 0.06696428571428571
$ws.Range("M15").Value = 0.004464285714285714
$ws.Range("N15").Value = 0.004464285714285714
$ws.Range("O15").Value = 0.0625
$ws.Range("S15").Value = 0.2321428571428572
$ws.Range("F16").Value = 0.008264462809917356
$ws.Range("H16").Value = 0.1652892561983471
$ws.Range("I16").Value = 0.07851239669421488
$ws.Range("J16").Value = 0.4049586776859504
$ws.Range("K16").Value = 0.1446280991735537
$ws.Range("M16").Value = 0.01652892561983471
$ws.Range("O16").Value = 0.04545454545454546
$ws.Range("S16").Value = 0.1363636363636364
$ws.Range("F17").Value = 0.01659751037344398
$ws.Range("H17").Value = 0.1390041493775934
$ws.Range("I17").Value = 0.1099585062240664
$ws.Range("J17").Value = 0.4211618257261411
$ws.Range("K17").Value = 0.0975103734439834
$ws.Range("M17").Value = 0.01452282157676349
$ws.Range("N17").Value = 0.002074688796680498
$ws.Range("O17").Value = 0.05601659751037345
$ws.Range("S17").Value = 0.1431535269709543
$ws.Range("F18").Value = 0.02409638554216868
$ws.Range("H18").Value = 0.1927710843373494
$ws.Range("I18").Value = 0.1204819277108434
$ws.Range("J18").Value = 0.3734939759036144
$ws.Range("K18").Value = 0.09036144578313253
$ws.Range("M18").Value = 0.006024096385542169
$ws.Range("O18").Value = 0.0783132530120482
$ws.Range("S18").Value = 0.1144578313253012
$ws.Range("F19").Value = 0.01536210680321873
$ws.Range("H19").Value = 0.1799561082662765
$ws.Range("I19").Value = 0.08558888076079005
$ws.Range("J19").Value = 0.3752743233357718
$ws.Range("K19").Value = 0.1177761521580102
$ws.Range("M19").Value = 0.01755669348939283
$ws.Range("N19").Value = 0.000731528895391368
$ws.Range("O19").Value = 0.06949524506217995
$ws.Range("S19").Value = 0.1382589612289686
